$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "0.9.20"
$wsMeta.Range("B8").Value = "2025-11-18T19:57:13-03:00"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AJ6").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`npositive-integer:Group size must be a positive integer (greater than 0) {`$this > 0}"
